$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..S of this export are numeric-looking trade figures ("0.01707000",
# "79699.33", ...) that the source system stores as plain TEXT (so leading/
# trailing zeros survive round-tripping). Plain `Range.Value = "79699.33"`
# lets Excel's smart-typing coerce that straight into a real number, which
# would lose the formatting. Forcing the cell to Text ("@") before the write
# keeps it a string; resetting `.Style` back to "Normal" afterwards drops the
# leftover quote-prefix/text formatting so the cell ends up with the same
# (default) style as its neighbours, matching the source file.
function Set-TextCell($row, $col, $text) {
    $rng = $ws.Cells.Item($row, $col)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# --- 1) Fix floating-point precision on the date_time serials of rows 164-167 ---
#     (sub-millisecond rounding differences in column A only; all other cells unchanged)
$ws.Cells.Item(164,1).Value = 45604.85270664352
$ws.Cells.Item(165,1).Value = 45604.85286658565
$ws.Cells.Item(166,1).Value = 45604.85302368055
$ws.Cells.Item(167,1).Value = 45604.85318097222

# --- 2) Append three new trade rows (168-170) pulled in from the small-portfolio bot run ---
#     Columns: A date_time | B asset | C transaction_type | D quantity | E price_per_asset |
#     F transaction_total | G reason | H profit_loss | I total_balance | J asset_balance |
#     K investment_usd | L returned_to_cash_usd | M cumulative_profit_loss | N usdt_balance |
#     O roi_percentage | P cumulative_profit_loss_percentage | Q market_value |
#     R decision_quality | S performance
$newRows = @(
    @{ Row = 168; DateSerial = 45606.58139018519; Values = @("BTCUSDT", "sell", "0.01707000", "79699.33", "1360.47", "Venda para lucro a curto prazo em carteira pequena com limite de vendas consecutivas", "1360.47", "69581.86", "0.17070000", "0.00", "1360.47", "1360.47", "69581.86", "0.00", "1.96", "0.00", "Good", "Profit of 1360.47") },
    @{ Row = 169; DateSerial = 45606.58153751157; Values = @("BTCUSDT", "sell", "0.01707000", "79697.13", "1360.43", "Venda para lucro a curto prazo em carteira pequena com limite de vendas consecutivas", "1360.43", "69581.86", "0.17070000", "0.00", "1360.43", "2720.90", "69581.86", "0.00", "3.91", "0.00", "Good", "Profit of 1360.43") },
    @{ Row = 170; DateSerial = 45606.58189127497; Values = @("BTCUSDT", "sell", "0.01365800", "79689.02", "1088.39", "Venda para lucro a curto prazo em carteira pequena com limite de vendas consecutivas", "1088.39", "72301.14", "0.13658000", "0.00", "1088.39", "1088.39", "72301.14", "0.00", "1.51", "0.00", "Good", "Profit of 1088.39") }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row

    # Column A is a real Excel date/time serial, formatted the same way as
    # the existing date_time column (rows 2-167).
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $dateCell.Value = $rowData.DateSerial

    $col = 2
    foreach ($val in $rowData.Values) {
        Set-TextCell $r $col $val
        $col++
    }
}

# Note: $ws.Cells.Item(<row>, 1) writes beyond the previous used range
# (rows stopped at 167), so Excel grows the sheet's dimension to A1:S170
# automatically once rows 168-170 are populated.
